# ---------------------------------------------------------------------------
# Re-applies the "sex" column addition to test_file sheet + trims Sheet2 back
# down to 3 rows, and restores Sheet2 as the active/selected sheet.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "test_file"
$ws2 = $wb.Worksheets.Item(2)   # "Sheet2"

# ---------------------------------------------------------------------------
# 1) Insert a new "sex" column between "age" (D) and "customer_type" (old E).
#    This shifts the existing customer_type column from E to F.
# ---------------------------------------------------------------------------
$ws1.Columns("E:E").Insert()

# ---------------------------------------------------------------------------
# 2) Fill in the header + the per-row sex values for the new column E.
# ---------------------------------------------------------------------------
$sexValues = @(
    "sex",
    "Male",
    "Other",
    "Female",
    "Male",
    "Female",
    "Male",
    "Female",
    "Male",
    "Female",
    "Male",
    "Male",
    "Male",
    "Female",
    "Female",
    "Female",
    "Female",
    "Female",
    "Male",
    "Female",
    "Male",
    "Male",
    "Female",
    "Female",
    "Female",
    "Male",
    "Male",
    "Male",
    "Male",
    "Female"
)

for ($i = 0; $i -lt $sexValues.Length; $i++) {
    $row = $i + 1
    $ws1.Cells.Item($row, 5).Value = $sexValues[$i]
}

# ---------------------------------------------------------------------------
# 3) Sheet2: remove the extra repeated rows 4-7, keeping only rows 1-3.
# ---------------------------------------------------------------------------
$ws2.Rows("4:7").Delete()

# ---------------------------------------------------------------------------
# 4) Make Sheet2 the active sheet/tab and select A3 on it (matches target).
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A3").Select()
